# SolicitudGrafica_CN_08_06_CO.xlsx
#
# The "Código guión o recurso" cell (C7) on the "Solicitud gráfica" sheet
# was holding the wrong script code (CN_08_01_CO_REC10, left over from a
# different guión/recurso). It must read the correct code for this
# guión: CN_08_06_CO. All the IMGxx file-name formulas in columns F/H
# already reference $C$7, so they recompute automatically once the cell
# is corrected (no image data changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Solicitud gráfica")
$ws.Activate()

$ws.Range("C7").Value = "CN_08_06_CO"

# Leave the selection on the cell that was edited.
[void]$ws.Range("C7").Select()
